$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 - this shifts the existing rows 9..89 down to 10..90
# and grows the used range to A1:R90 (matching the dimension change in the diff).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C9").Value = "Los Lagos"
$ws.Range("D9").Value = 44761
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 100112031
$ws.Range("G9").Value = "Poroto verde"
$ws.Range("H9").Value = "Magnum"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 37000
$ws.Range("L9").Value = 37000
$ws.Range("M9").Value = 37000
$ws.Range("N9").Value = "`$/malla 25 kilos"
$ws.Range("O9").Value = "Perú"
$ws.Range("P9").Value = 1480
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
